$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 01:52"

# Direct numeric updates (no row/label shift)
$ws.Range("B4").Value = 817187
$ws.Range("C4").Value = 24428
$ws.Range("D4").Value = 82860
$ws.Range("E4").Value = 689098
$ws.Range("G4").Value = 2715
$ws.Range("H4").Value = 45229

$ws.Range("B8").Value = 148453
$ws.Range("C8").Value = 1388
$ws.Range("E8").Value = 48167
$ws.Range("G8").Value = 224
$ws.Range("H8").Value = 5086

$ws.Range("D16").Value = 13188
$ws.Range("E16").Value = 23400
$ws.Range("G16").Value = 144
$ws.Range("H16").Value = 1834

$ws.Range("B88").Value = 901
$ws.Range("C88").Value = 17
$ws.Range("D88").Value = 170
$ws.Range("E88").Value = 693
$ws.Range("F88").Value = 25

$ws.Range("B102").Value = 543
$ws.Range("C102").Value = 8
$ws.Range("D102").Value = 324
$ws.Range("E102").Value = 207
$ws.Range("F102").Value = 10
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 12

$ws.Range("B160").Value = 57
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 35
$ws.Range("E160").Value = 22

$ws.Range("B166").Value = 42
$ws.Range("C166").Value = 11
$ws.Range("E166").Value = 38

# Nigeria moves up in ranking: rows 90-96 shift (label + data) down by one,
# row 90 gets fresh Nigeria data; row 97 (Niger) unaffected
$ws.Range("A90").Value = "Nigeria"
$ws.Range("B90").Value = 782
$ws.Range("C90").Value = 117
$ws.Range("D90").Value = 197
$ws.Range("E90").Value = 560
$ws.Range("F90").Value = 2
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 25

$ws.Range("A91").Value = "Letonia"
$ws.Range("B91").Value = 748
$ws.Range("C91").Value = 9
$ws.Range("D91").Value = 133
$ws.Range("E91").Value = 606
$ws.Range("F91").Value = 3
$ws.Range("G91").Value = 4
$ws.Range("H91").Value = 9

$ws.Range("A92").Value = "Principado de Andorra"
$ws.Range("B92").Value = 717
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 282
$ws.Range("E92").Value = 398
$ws.Range("F92").Value = 17
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 37

$ws.Range("A93").Value = "Crucero"
$ws.Range("B93").Value = 712
$ws.Range("C93").Value = 0
$ws.Range("D93").Value = 644
$ws.Range("E93").Value = 55
$ws.Range("F93").Value = 7
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 13

$ws.Range("A94").Value = "Guinea"
$ws.Range("B94").Value = 688
$ws.Range("C94").Value = 66
$ws.Range("D94").Value = 127
$ws.Range("E94").Value = 555
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 6

$ws.Range("A95").Value = "Libano"
$ws.Range("B95").Value = 677
$ws.Range("C95").Value = 0
$ws.Range("D95").Value = 108
$ws.Range("E95").Value = 548
$ws.Range("F95").Value = 27
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 21

$ws.Range("A96").Value = "Costa Rica"
$ws.Range("B96").Value = 669
$ws.Range("C96").Value = 7
$ws.Range("D96").Value = 150
$ws.Range("E96").Value = 513
$ws.Range("F96").Value = 6
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 6

# Sudan del Sur / Santo Tome y Principe swap (rank tie, data unchanged)
$ws.Range("A212").Value = "Sudan del Sur"
$ws.Range("A213").Value = "Santo Tome y Principe"

# San Pedro y Miquelon / Yemen swap (rank tie, data unchanged)
$ws.Range("A215").Value = "San Pedro y Miquelon"
$ws.Range("A216").Value = "Yemen"
